# Apply the per-row Coin/Link/Price/Volume(1h) updates captured in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose text looks like a plain decimal number get an apostrophe prefix so
# Excel stores them as text (matching the source inline-string cells) instead of
# silently re-typing them as the Number data type.
$textForcedCells = New-Object System.Collections.Generic.List[string]

$ws.Range("D2").Value = "30.486.66"
$ws.Range("E2").Value = "  -1.06%  "
$ws.Range("D3").Value = "1.912.28"
$ws.Range("E3").Value = "  -1.54%  "
$ws.Range("D4").Value = "'1.001"; [void]$textForcedCells.Add("D4")
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'239.33"; [void]$textForcedCells.Add("D5")
$ws.Range("E5").Value = "  -1.49%  "
$ws.Range("D6").Value = "'1.000"; [void]$textForcedCells.Add("D6")
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "'0.4783"; [void]$textForcedCells.Add("D7")
$ws.Range("D8").Value = "'0.2842"; [void]$textForcedCells.Add("D8")
$ws.Range("E8").Value = "  -3.83%  "
$ws.Range("D9").Value = "'0.06711"; [void]$textForcedCells.Add("D9")
$ws.Range("E9").Value = "  -2.44%  "
$ws.Range("D10").Value = "'18.82"; [void]$textForcedCells.Add("D10")
$ws.Range("E10").Value = "  -3.17%  "
$ws.Range("D11").Value = "'102.11"; [void]$textForcedCells.Add("D11")
$ws.Range("E11").Value = "  -3.96%  "
$ws.Range("D12").Value = "'0.07706"; [void]$textForcedCells.Add("D12")
$ws.Range("E12").Value = "  -0.32%  "
$ws.Range("D13").Value = "1.915.23"
$ws.Range("E13").Value = "  -0.70%  "
$ws.Range("D14").Value = "'5.207"; [void]$textForcedCells.Add("D14")
$ws.Range("E14").Value = "  -2.77%  "
$ws.Range("D15").Value = "'0.6711"; [void]$textForcedCells.Add("D15")
$ws.Range("E15").Value = "  -3.85%  "
$ws.Range("D16").Value = "'265.00"; [void]$textForcedCells.Add("D16")
$ws.Range("E16").Value = "  -4.32%  "
$ws.Range("D17").Value = "30.515.65"
$ws.Range("E17").Value = "  -0.98%  "
$ws.Range("D18").Value = "'1.001"; [void]$textForcedCells.Add("D18")
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D19").Value = "'0.000007474"; [void]$textForcedCells.Add("D19")
$ws.Range("E19").Value = "  -3.23%  "
$ws.Range("D20").Value = "'12.68"; [void]$textForcedCells.Add("D20")
$ws.Range("E20").Value = "  -3.34%  "
$ws.Range("D21").Value = "'5.410"; [void]$textForcedCells.Add("D21")
$ws.Range("E21").Value = "  -1.10%  "
$ws.Range("D22").Value = "'1.001"; [void]$textForcedCells.Add("D22")
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").Value = "'6.301"; [void]$textForcedCells.Add("D23")
$ws.Range("E23").Value = "  -3.57%  "
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").Value = "'9.373"; [void]$textForcedCells.Add("D24")
$ws.Range("E24").Value = "  -3.67%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "'167.16"; [void]$textForcedCells.Add("D25")
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "'19.19"; [void]$textForcedCells.Add("D26")
$ws.Range("E26").Value = "  -2.39%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").Value = "'2.064"; [void]$textForcedCells.Add("D27")
$ws.Range("E27").Value = "  -4.86%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "'1.384"; [void]$textForcedCells.Add("D28")
$ws.Range("E28").Value = "  -0.49%  "
$ws.Range("B29").Value = "Stellar"
$ws.Range("C29").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D29").Value = "'0.1002"; [void]$textForcedCells.Add("D29")
$ws.Range("E29").Value = "  -4.39%  "
$ws.Range("D30").Value = "'4.620"; [void]$textForcedCells.Add("D30")
$ws.Range("E30").Value = "  +1.08%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.511"; [void]$textForcedCells.Add("D31")
$ws.Range("E31").Value = "  -2.84%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "'4.227"; [void]$textForcedCells.Add("D32")
$ws.Range("E32").Value = "  -3.40%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.04730"; [void]$textForcedCells.Add("D33")
$ws.Range("E33").Value = "  -2.62%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "'0.7266"; [void]$textForcedCells.Add("D34")
$ws.Range("E34").Value = "  -3.61%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'1.108"; [void]$textForcedCells.Add("D35")
$ws.Range("E35").Value = "  -4.42%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "'2.717"; [void]$textForcedCells.Add("D36")
$ws.Range("E36").Value = "  -0.55%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.01916"; [void]$textForcedCells.Add("D37")
$ws.Range("E37").Value = "  -4.01%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").Value = "'2.612"; [void]$textForcedCells.Add("D38")
$ws.Range("E38").Value = "  -1.67%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "'6.268"; [void]$textForcedCells.Add("D39")
$ws.Range("E39").Value = "  -3.92%  "
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").Value = "'75.10"; [void]$textForcedCells.Add("D40")
$ws.Range("E40").Value = "  -3.54%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "'1.969"; [void]$textForcedCells.Add("D41")
$ws.Range("E41").Value = "  -6.30%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'0.8606"; [void]$textForcedCells.Add("D42")
$ws.Range("E42").Value = "  -5.21%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "'105.01"; [void]$textForcedCells.Add("D43")
$ws.Range("E43").Value = "  -3.03%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "'0.4265"; [void]$textForcedCells.Add("D44")
$ws.Range("E44").Value = "  -3.01%  "
$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").Value = "'0.9998"; [void]$textForcedCells.Add("D45")
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").Value = "'7.412"; [void]$textForcedCells.Add("D46")
$ws.Range("E46").Value = "  -4.71%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "'927.62"; [void]$textForcedCells.Add("D47")
$ws.Range("E47").Value = "  -7.20%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Value = "'0.1201"; [void]$textForcedCells.Add("D48")
$ws.Range("E48").Value = "  -3.68%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "'34.78"; [void]$textForcedCells.Add("D49")
$ws.Range("E49").Value = "  -3.40%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'8.774"; [void]$textForcedCells.Add("D50")
$ws.Range("E50").Value = "  -5.40%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.05754"; [void]$textForcedCells.Add("D51")
$ws.Range("E51").Value = "  +0.40%  "

# Drop the implicit "Text" number-format style the apostrophe prefix applies so the
# cells end up with no style attribute, same as the original inline-string cells.
foreach ($cellRef in $textForcedCells) {
    $ws.Range($cellRef).Style = "Normal"
}
